# Updates Price (D) and Volume(1h) (E) columns on the cryptos sheet
# to reflect the latest scrape, per the automated GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.026.29'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '2.504.06'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '534.33'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.77'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '2.515.97'
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  -2.60%  '
$ws.Range('D14').Value = '2.963.88'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.13'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = '58.960.58'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '2.514.53'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.01'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '324.35'
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.83'
$ws.Range('E24').Value = '  +5.55%  '
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -2.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.75'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '167.85'
$ws.Range('E32').Value = '  +4.55%  '
$ws.Range('E33').Value = '  +4.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -4.30%  '
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.78'
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.825'
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '279.42'
$ws.Range('E43').Value = '  -2.12%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.603'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.87'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '128.09'
$ws.Range('E47').Value = '  +3.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0928'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0222'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.32'
$ws.Range('E51').Value = '  -1.25%  '
